# Append OCR'd pages-viewed rows (140-144) to Sheet1, matching the commit's
# "ocr updates" data. Columns: A=original paper, B=figure name,
# C=figure number, D=year, E=page number, F=image rotation, G=is_viewed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 140; A = 'Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.pdf'; B = 'F2_P12_Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.png'; C = 2; D = 2021; E = 13; F = 0; G = $false },
    @{ r = 141; A = 'Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.pdf'; B = 'F2_P11_Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.png'; C = 2; D = 2021; E = 12; F = 0; G = $false },
    @{ r = 142; A = 'Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.pdf'; B = 'F3_P12_Allard-Poesi & Hollet-Haudebert_2021_HR_The Sound of Silence Measuring Suffering at Work.png'; C = 3; D = 2021; E = 13; F = 0; G = $false },
    @{ r = 143; A = 'Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.pdf'; B = 'F1_P6_Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.png'; C = 1; D = 2021; E = 7; F = 0; G = $false },
    @{ r = 144; A = 'Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.pdf'; B = 'F2_P6_Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.png'; C = 2; D = 2021; E = 7; F = 0; G = $false }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
